# Updated cryptos list: applies new Price (D) and Volume(1h) (E) values
# to rows 2-51 on the active worksheet, matching the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.122.14'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '2.614.22'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''589.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").Value = '''164.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.85%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''0.530'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = '2.613.77'
$ws.Range("E9").Value = '  -1.11%  '
$ws.Range("E10").Value = '  -5.11%  '
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").Value = '''5.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").Value = '''27.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.04%  '
$ws.Range("D15").Value = '3.089.38'
$ws.Range("E15").Value = '  -0.94%  '
$ws.Range("D16").Value = '''0.0000179'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.88%  '
$ws.Range("D17").Value = '67.033.36'
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").Value = '2.605.73'
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").Value = '''11.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.32%  '
$ws.Range("D20").Value = '''7.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.65%  '
$ws.Range("D21").Value = '''355.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.36%  '
$ws.Range("E22").Value = '  -3.09%  '
$ws.Range("D23").Value = '''4.63'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.42%  '
$ws.Range("D24").Value = '''10.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.44%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '''1.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.17%  '
$ws.Range("D27").Value = '''69.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.24%  '
$ws.Range("D28").Value = '2.748.39'
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").Value = '0.0₃0993'
$ws.Range("E30").Value = '  -3.34%  '
$ws.Range("D31").Value = '''543.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.73%  '
$ws.Range("D32").Value = '''7.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.17%  '
$ws.Range("D33").Value = '''1.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.40%  '
$ws.Range("E34").Value = '  -3.07%  '
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").Value = '''1.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.70%  '
$ws.Range("D38").Value = '''157.74'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("D39").Value = '''18.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.71%  '
$ws.Range("D40").Value = '''0.363'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.64%  '
$ws.Range("E41").Value = '  +1.71%  '
$ws.Range("E42").Value = '  -2.09%  '
$ws.Range("D43").Value = '''5.12'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.33%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '''2.40'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.19%  '
$ws.Range("D46").Value = '0.0₆0296'
$ws.Range("E46").Value = '  -1.42%  '
$ws.Range("D47").Value = '''151.02'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("D48").Value = '''0.574'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.63%  '
$ws.Range("E49").Value = '  -3.13%  '
$ws.Range("D50").Value = '''1.70'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.84%  '
$ws.Range("E51").Value = '  -1.43%  '
